$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '70.284.14'
$ws.Range('E2').Value = '  -2.69%  '
$ws.Range('D3').Value = '2.522.82'
$ws.Range('E3').Value = '  -4.83%  '
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '575.56'
$ws.Range('E5').Value = '  -3.49%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '169.63'
$ws.Range('E6').Value = '  -2.77%  '
$ws.Range('E8').Value = '  -2.51%  '
$ws.Range('D9').Value = '2.522.60'
$ws.Range('E9').Value = '  -4.77%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.160'
$ws.Range('E10').Value = '  -5.40%  '
$ws.Range('E11').Value = '  -0.68%  '
$ws.Range('E12').Value = '  -3.69%  '
$ws.Range('E13').Value = '  -3.23%  '
$ws.Range('D14').Value = '2.984.03'
$ws.Range('E14').Value = '  -4.97%  '
$ws.Range('D15').Value = '70.160.28'
$ws.Range('E15').Value = '  -2.73%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000178'
$ws.Range('E16').Value = '  -3.31%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '25.07'
$ws.Range('E17').Value = '  -4.47%  '
$ws.Range('D18').Value = '2.524.10'
$ws.Range('E18').Value = '  -4.80%  '
$ws.Range('E19').Value = '  -5.46%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.61'
$ws.Range('E20').Value = '  -6.30%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '355.58'
$ws.Range('E21').Value = '  -3.86%  '
$ws.Range('E22').Value = '  -5.26%  '
$ws.Range('E23').Value = '  -2.60%  '
$ws.Range('E24').Value = '  +0.16%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '69.06'
$ws.Range('E25').Value = '  -4.21%  '
$ws.Range('E26').Value = '  -5.06%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.20'
$ws.Range('E27').Value = '  -5.44%  '
$ws.Range('D28').Value = '2.651.50'
$ws.Range('E28').Value = '  -5.16%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.00'
$ws.Range('E29').Value = '  +0.04%  '
$ws.Range('D30').Value = '0.0₃0912'
$ws.Range('E30').Value = '  -5.74%  '
$ws.Range('E31').Value = '  -2.90%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '483.69'
$ws.Range('E32').Value = '  -2.97%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.31'
$ws.Range('E33').Value = '  +1.29%  '
$ws.Range('E34').Value = '  -3.25%  '
$ws.Range('E35').Value = '  -0.08%  '
$ws.Range('E36').Value = '  +4.91%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '155.49'
$ws.Range('E37').Value = '  -4.50%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '18.88'
$ws.Range('E38').Value = '  -0.20%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '18.60'
$ws.Range('E39').Value = '  -4.48%  '
$ws.Range('B41').Value = 'Stacks'
$ws.Range('C41').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.65'
$ws.Range('E41').Value = '  -6.47%  '
$ws.Range('B42').Value = 'PolygonEcosystemToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.321'
$ws.Range('E42').Value = '  -3.01%  '
$ws.Range('E43').Value = '  -4.77%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.22'
$ws.Range('E44').Value = '  -11.28%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.38'
$ws.Range('E45').Value = '  -7.49%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '38.27'
$ws.Range('E46').Value = '  -2.99%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '142.59'
$ws.Range('E47').Value = '  -8.37%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.53'
$ws.Range('E48').Value = '  -5.50%  '
$ws.Range('E49').Value = '  -5.25%  '
$ws.Range('E50').Value = '  -5.99%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.600'
$ws.Range('E51').Value = '  -0.57%  '
